# Replace the PSSM score matrix (B2:K21) with updated supplemental-figure values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    @(-18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833),
    @(-18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, 1.286362968724222, -18.86169274946833, -18.86169274946833),
    @(-18.86169274946833, -18.86169274946833, 2.031899861687916, -18.86169274946833, 3.390326126951632, -18.86169274946833, 1.476522415861057, -18.86169274946833, 0.9460976652014469, -18.86169274946833),
    @(-18.86169274946833, 4.321925217732153, -18.86169274946833, -18.86169274946833, -18.86169274946833, 2.812856868258528, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833),
    @(-18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833),
    @(2.460202874729149, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833),
    @(-18.86169274946833, -18.86169274946833, -18.86169274946833, 1.779132042259539, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833),
    @(3.857684581436429, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833),
    @(-18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, 1.764498421248164, -18.86169274946833, 2.245228596672643),
    @(-18.86169274946833, -18.86169274946833, -18.86169274946833, 2.894293799938684, -18.86169274946833, 2.854505849879029, -18.86169274946833, -18.86169274946833, -18.86169274946833, 1.952480527253549),
    @(-18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833),
    @(-18.86169274946833, -18.86169274946833, -18.86169274946833, 2.52993736522065, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, 1.659831972847947, 1.725191924148898),
    @(-18.86169274946833, -18.86169274946833, 1.358058353656082, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, 1.944193578007145),
    @(-18.86169274946833, -18.86169274946833, 1.507083314013883, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833),
    @(-18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, 1.908417000490539, -18.86169274946833),
    @(-18.86169274946833, -18.86169274946833, 1.764409278561828, -18.86169274946833, -18.86169274946833, -18.86169274946833, 2.075736390369192, 2.092024148980915, 2.520340692264334, -18.86169274946833),
    @(-18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, -18.86169274946833, 1.998622094800338, 2.030062034305842, 2.43892402184632, -18.86169274946833),
    @(-18.86169274946833, -18.86169274946833, 2.140980135088863, -18.86169274946833, -18.86169274946833, -18.86169274946833, 1.642533668566871, 1.794684102785578, -18.86169274946833, -18.86169274946833),
    @(-18.86169274946833, -18.86169274946833, 1.431263527024221, -18.86169274946833, 3.250119496721957, -18.86169274946833, 1.665590474594026, 1.230332671102149, -18.86169274946833, 2.082016090562137),
    @(-18.86169274946833, -18.86169274946833, -18.86169274946833, 1.747389474330971, -18.86169274946833, 2.521216718939914, 1.438745838353981, -18.86169274946833, -18.86169274946833, -18.86169274946833)
)

$arr = New-Object 'object[,]' 20,10
for ($r = 0; $r -lt $newValues.Count; $r++) {
    for ($c = 0; $c -lt $newValues[$r].Count; $c++) {
        $arr[$r, $c] = $newValues[$r][$c]
    }
}

$ws.Range("B2:K21").Value = $arr
